# Add team record (Wins/Losses/Ties) columns to the data sheet.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Header row (row 1): new columns AD, AE, AF
$ws.Range("AD1").Value = "Wins"
$ws.Range("AE1").Value = "Losses"
$ws.Range("AF1").Value = "Ties"

# Copy the header style from an existing header cell (e.g. AC1) so the
# new header cells match the bold/bordered/centered style used for the
# rest of row 1.
$ws.Range("AC1").Copy()
$ws.Range("AD1:AF1").PasteSpecial(-4122) # xlPasteFormats
$excel.CutCopyMode = 0

# Data rows 2-51: constant team record values for every player row.
$lastRow = 51
for ($r = 2; $r -le $lastRow; $r++) {
    $ws.Cells.Item($r, 30).Value = 72  # AD = column 30 -> Wins
    $ws.Cells.Item($r, 31).Value = 90  # AE = column 31 -> Losses
    $ws.Cells.Item($r, 32).Value = 0   # AF = column 32 -> Ties
}
